$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.527.97"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "'3.844.80"
$ws.Range("E3").Value = "  -2.38%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'518.41"
$ws.Range("E5").Value = "  +5.09%  "
$ws.Range("D6").Value = "'140.37"
$ws.Range("E6").Value = "  -4.85%  "
$ws.Range("D7").Value = "'0.606"
$ws.Range("E7").Value = "  -2.65%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'0.710"
$ws.Range("E9").Value = "  -3.20%  "
$ws.Range("D10").Value = "'0.167"
$ws.Range("E10").Value = "  -5.30%  "
$ws.Range("D11").Value = "'0.0000319"
$ws.Range("E11").Value = "  -9.21%  "
$ws.Range("D12").Value = "'41.56"
$ws.Range("E12").Value = "  -4.27%  "
$ws.Range("D13").Value = "'10.29"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").Value = "'4.459.19"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").Value = "'21.25"
$ws.Range("E15").Value = "  +6.89%  "
$ws.Range("D16").Value = "'3.868.41"
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("D17").Value = "'13.97"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").Value = "'68.445.04"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").Value = "'412.83"
$ws.Range("E21").Value = "  -6.32%  "
$ws.Range("D22").Value = "'3.45"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'13.97"
$ws.Range("E23").Value = "  -3.70%  "
$ws.Range("D24").Value = "'12.03"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "'86.44"
$ws.Range("E25").Value = "  -2.63%  "
$ws.Range("D26").Value = "'4.00"
$ws.Range("E26").Value = "  +5.86%  "
$ws.Range("D27").Value = "'10.38"
$ws.Range("E27").Value = "  -7.06%  "
$ws.Range("D28").Value = "'35.28"
$ws.Range("E28").Value = "  -4.96%  "
$ws.Range("D29").Value = "'13.28"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").Value = "'675.36"
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("D31").Value = "'7.00"
$ws.Range("E31").Value = "  +15.39%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'2.84"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.124"
$ws.Range("E33").Value = "  -4.88%  "
$ws.Range("D34").Value = "'66.76"
$ws.Range("E34").Value = "  +8.29%  "
$ws.Range("D35").Value = "'0.447"
$ws.Range("E35").Value = "  -5.50%  "
$ws.Range("D36").Value = "'0.0₃0842"
$ws.Range("E36").Value = "  -7.56%  "
$ws.Range("D37").Value = "'39.34"
$ws.Range("E37").Value = "  -3.30%  "
$ws.Range("D38").Value = "'3.40"
$ws.Range("E38").Value = "  +11.06%  "
$ws.Range("D39").Value = "'0.147"
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "'2.87"
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("D43").Value = "'0.0473"
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("D44").Value = "'3.14"
$ws.Range("E44").Value = "  +4.91%  "
$ws.Range("D45").Value = "'3.41"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").Value = "'0.000277"
$ws.Range("E47").Value = "  +14.46%  "
$ws.Range("D48").Value = "'3.00"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").Value = "'3.28"
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'142.65"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'8.71"
$ws.Range("E51").Value = "  +2.78%  "
